# Fruta / hortaliza, semanal
# A new weekly price record (Terminal Hortofrutícola Agro Chillán - Papa,
# "1a (guarda lavada)", Región de Los Lagos) is inserted as row 141,
# pushing the existing rows 141-209 down to 142-210.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new, blank row above the current row 141 - everything below
# (old rows 141..209) shifts down to 142..210.
$ws.Rows(141).Insert()

$row = 141
$ws.Cells.Item($row, 1).Value  = 7
$ws.Cells.Item($row, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item($row, 3).Value  = "Ñuble"
$ws.Cells.Item($row, 4).Value  = 44489
$ws.Cells.Item($row, 5).Value  = 16
$ws.Cells.Item($row, 6).Value  = 100114001
$ws.Cells.Item($row, 7).Value  = "Papa"
$ws.Cells.Item($row, 8).Value  = "Patagonia"
$ws.Cells.Item($row, 9).Value  = "1a (guarda lavada)"
$ws.Cells.Item($row, 10).Value = 160
$ws.Cells.Item($row, 11).Value = 9500
$ws.Cells.Item($row, 12).Value = 10000
$ws.Cells.Item($row, 13).Value = 9750
$ws.Cells.Item($row, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item($row, 15).Value = "Región de Los Lagos"
$ws.Cells.Item($row, 16).Value = 390
$ws.Cells.Item($row, 17).Value = 25
$ws.Cells.Item($row, 18).Value = "Hortaliza"
